$d = $word.ActiveDocument

# The four inline pictures (two "PearsonLogo.png" images living in the
# footers, two "BTec_Logo-Orange" images living in the headers) each carry
# a display "name" in both the <wp:docPr> and the matching <pic:cNvPr>
# element. InlineShape objects don't expose a settable .Name property in
# the Word object model, so we rewrite those two attributes directly via
# the document's WordOpenXML (round-trips losslessly) using unique,
# id-anchored search/replace strings - one pair per picture.

$xml = $d.WordOpenXML

# Footer picture with id="2": image1.png -> image2.png
$old1 = '<wp:docPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="2" name="image1.png"/><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="0" name="image1.png"/>'
$new1 = '<wp:docPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="2" name="image2.png"/><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="0" name="image2.png"/>'
$xml = $xml.Replace($old1, $new1)

# Footer picture with id="4": image1.png -> image2.png
$old2 = '<wp:docPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="4" name="image1.png"/><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="0" name="image1.png"/>'
$new2 = '<wp:docPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="4" name="image2.png"/><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="0" name="image2.png"/>'
$xml = $xml.Replace($old2, $new2)

# Header picture with id="1": image2.jpg -> image1.jpg
$old3 = '<wp:docPr descr="BTec_Logo-Orange" id="1" name="image2.jpg"/><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr descr="BTec_Logo-Orange" id="0" name="image2.jpg"/>'
$new3 = '<wp:docPr descr="BTec_Logo-Orange" id="1" name="image1.jpg"/><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr descr="BTec_Logo-Orange" id="0" name="image1.jpg"/>'
$xml = $xml.Replace($old3, $new3)

# Header picture with id="3": image2.jpg -> image1.jpg
$old4 = '<wp:docPr descr="BTec_Logo-Orange" id="3" name="image2.jpg"/><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr descr="BTec_Logo-Orange" id="0" name="image2.jpg"/>'
$new4 = '<wp:docPr descr="BTec_Logo-Orange" id="3" name="image1.jpg"/><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr descr="BTec_Logo-Orange" id="0" name="image1.jpg"/>'
$xml = $xml.Replace($old4, $new4)

$d.WordOpenXML = $xml
